# Update the "想去人数" (want-to-go count) values on the "展览" and
# "全部类型" sheets. Row 3 (HP国风动漫游戏嘉年华) goes 185 -> 186 and
# row 4 (动漫游戏展) goes 137 -> 138 on both sheets.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 186
    $ws.Range("F4").Value = 138
}
